$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 262, pushing existing rows 262..302 down to 263..303.
$ws.Rows.Item(262).Insert()

# Populate the new row 262 with a new weekly price record (clone of the
# surrounding record's fixed fields, with updated date/volume/price fields).
$ws.Cells.Item(262, 1).Value = 5
$ws.Cells.Item(262, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(262, 3).Value = "Maule"
$ws.Cells.Item(262, 4).Value = 45077
$ws.Cells.Item(262, 5).Value = 7
$ws.Cells.Item(262, 6).Value = 100112017
$ws.Cells.Item(262, 7).Value = "Apio"
$ws.Cells.Item(262, 8).Value = "Americana (o)"
$ws.Cells.Item(262, 9).Value = "Primera"
$ws.Cells.Item(262, 10).Value = 700
$ws.Cells.Item(262, 11).Value = 5500
$ws.Cells.Item(262, 12).Value = 5500
$ws.Cells.Item(262, 13).Value = 5500
$ws.Cells.Item(262, 14).Value = "$/docena de matas"
$ws.Cells.Item(262, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(262, 16).Value = 917
$ws.Cells.Item(262, 17).Value = 6
$ws.Cells.Item(262, 18).Value = "Hortaliza"
